$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Activate()

# Delete row 74 (the row for DOI 10.3389/fmicb.2021.696522), shifting all
# subsequent rows up by one. This removes its shared-string entries too and
# re-numbers the worksheet's used range from A1:K132 down to A1:K131.
$ws.Rows.Item(74).Delete()

# Match the saved selection/scroll state: the active cell sits on the (now
# shifted) row that used to be row 75, selected as a full row.
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("A74:XFD74").Select()
